$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string rich-text cells collapse to plain text,
# which is what Excel does when a rich-text cell's formula-bar value is
# retyped and every run shares identical formatting, as is the case here).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# ---------------------------------------------------------------------------
# Row 14 - Murder (value-only change)
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -87.5

# ---------------------------------------------------------------------------
# Row 15 - Rape (C15/D15/E15 flip from text placeholders to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Copy()
$ws.Range("C15:D15").PasteSpecial(-4122)
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -20
$ws.Range("I15").Value = 33
$ws.Range("J15").Value = 33
$ws.Range("L15").Value = 13.793103448275
$ws.Range("M15").Value = 57.142857142857
$ws.Range("N15").Value = 32

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 337
$ws.Range("J16").Value = 425
$ws.Range("K16").Value = -20.705882352941
$ws.Range("L16").Value = -10.133333333333
$ws.Range("M16").Value = 17.013888888888
$ws.Range("N16").Value = -72.888173773129

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -57.894736842105
$ws.Range("F17").Value = 58
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = -15.942028985507
$ws.Range("I17").Value = 487
$ws.Range("J17").Value = 676
$ws.Range("K17").Value = -27.958579881656
$ws.Range("L17").Value = -18.833333333333
$ws.Range("M17").Value = 145.959595959596
$ws.Range("N17").Value = 24.552429667519

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = -35.714285714285
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 199
$ws.Range("K18").Value = -33.165829145728
$ws.Range("L18").Value = -13.636363636363
$ws.Range("M18").Value = -42.173913043478
$ws.Range("N18").Value = -92.4

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 93
$ws.Range("H19").Value = -51.612903225806
$ws.Range("I19").Value = 616
$ws.Range("J19").Value = 918
$ws.Range("K19").Value = -32.897603485838
$ws.Range("L19").Value = -28.621089223638
$ws.Range("M19").Value = 40.639269406392
$ws.Range("N19").Value = -37.398373983739

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 196
$ws.Range("J20").Value = 235
$ws.Range("K20").Value = -16.595744680851
$ws.Range("L20").Value = -28.985507246376
$ws.Range("M20").Value = 60.655737704918
$ws.Range("N20").Value = -89.201101928374

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 63
$ws.Range("E21").Value = -36.507936507936
$ws.Range("F21").Value = 172
$ws.Range("G21").Value = 237
$ws.Range("H21").Value = -27.426160337552
$ws.Range("I21").Value = 1805
$ws.Range("J21").Value = 2488
$ws.Range("K21").Value = -27.451768488746
$ws.Range("L21").Value = -21.521739130434
$ws.Range("M21").Value = 38.846153846153
$ws.Range("N21").Value = -71.036585365853

# ---------------------------------------------------------------------------
# Row 22 - Transit (C22 flips number -> text; D22/E22 flip text -> number)
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("H22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C22").Value = "0"
$ws.Range("N22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("J22").Value = 41
$ws.Range("K22").Value = 14.634146341463
$ws.Range("L22").Value = 9.302325581395

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 53
$ws.Range("E24").Value = -33.962264150943
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 178
$ws.Range("H24").Value = -14.606741573033
$ws.Range("I24").Value = 1907
$ws.Range("J24").Value = 2463
$ws.Range("K24").Value = -22.574096630125
$ws.Range("L24").Value = -15.917107583774
$ws.Range("M24").Value = 34.012649332396

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 41
$ws.Range("E25").Value = -58.536585365853
$ws.Range("F25").Value = 66
$ws.Range("G25").Value = 121
$ws.Range("H25").Value = -45.454545454545
$ws.Range("I25").Value = 1309
$ws.Range("J25").Value = 1931
$ws.Range("K25").Value = -32.211289487312
$ws.Range("L25").Value = -21.569802276812

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 27
$ws.Range("E26").Value = -74.074074074074
$ws.Range("F26").Value = 77
$ws.Range("G26").Value = 100
$ws.Range("H26").Value = -23
$ws.Range("I26").Value = 850
$ws.Range("J26").Value = 1108
$ws.Range("K26").Value = -23.285198555956
$ws.Range("L26").Value = -8.207343412527
$ws.Range("M26").Value = 62.523900573613

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape* (D27/E27 flip text -> number)
# ---------------------------------------------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -44.444444444444
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -12
$ws.Range("L27").Value = -2.222222222222

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -26.666666666666
$ws.Range("I28").Value = 98
$ws.Range("J28").Value = 127
$ws.Range("K28").Value = -22.834645669291
$ws.Range("L28").Value = -14.035087719298

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic. (value-only change)
# ---------------------------------------------------------------------------
$ws.Range("N29").Value = -84

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc. (value-only change)
# ---------------------------------------------------------------------------
$ws.Range("N30").Value = -89.130434782608

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes (G31/H31 flip number -> text placeholders)
# ---------------------------------------------------------------------------
$ws.Range("G31").Value = "'0"
$ws.Range("H31").Value = "***.*"
$ws.Range("C31").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("D31").Copy()
$ws.Range("H31").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 33 - Traffic Fatalities (D33/E33 flip number -> text placeholders)
# ---------------------------------------------------------------------------
$ws.Range("D33").Value = "'0"
$ws.Range("E33").Value = "***.*"
$ws.Range("C33").Copy()
$ws.Range("D33:E33").PasteSpecial(-4122)
